$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a couple of mis-entered ratio values
$ws.Range("W27").Value2 = 0.05535055350553506
$ws.Range("W28").Value2 = 0.05535055350553506

# Rename the "C_conversion_ratio" header to "C_ratio"
$ws.Range("AL1").Value2 = "C_ratio"

# Rows 149-159 actually belong one row higher (148-158); shift that whole
# block up by one row (columns B..AL only - column A is an independent
# running counter and must stay where it is), then drop the now-duplicated
# last row (159).
$src = $ws.Range("B149:AL159")
$dst = $ws.Range("B148")
$src.Copy($dst)

$ws.Rows.Item(159).Delete()
